# "Added data provider and reporting"
# Renames the sheet to userRegistration, fills in a user-registration
# data table (6 columns x 3 rows incl. header), adds mailto hyperlinks
# on the two email cells, sets column widths, and restores the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet ---
$ws.Name = "userRegistration"

# --- Header row ---
$ws.Cells.Item(1,1).Value = "firstName"
$ws.Cells.Item(1,2).Value = "lastName"
$ws.Cells.Item(1,3).Value = "email"
$ws.Cells.Item(1,4).Value = "phone"
$ws.Cells.Item(1,5).Value = "password"
$ws.Cells.Item(1,6).Value = "confirmPwd"

# --- Row 2: Prem Mishra ---
$ws.Cells.Item(2,1).Value = "Prem"
$ws.Cells.Item(2,2).Value = "Mishra"
$ws.Cells.Item(2,3).Value = "Prem@gmail.com"
$ws.Cells.Item(2,4).Value = 2234456567
$ws.Cells.Item(2,5).Value = "prem123"
$ws.Cells.Item(2,6).Value = "prem123"

# --- Row 3: Rajan Sharma ---
$ws.Cells.Item(3,1).Value = "Rajan"
$ws.Cells.Item(3,2).Value = "Sharma"
$ws.Cells.Item(3,3).Value = "Rajan@gmail.com"
$ws.Cells.Item(3,4).Value = 9678765321
$ws.Cells.Item(3,5).Value = "rajan123"
$ws.Cells.Item(3,6).Value = "rajan123"

# --- Hyperlinks on the email cells (mailto:) ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Prem@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Rajan@gmail.com")

# --- Column widths to fit the new data ---
$ws.Columns.Item(1).ColumnWidth = 13.76
$ws.Columns.Item(2).ColumnWidth = 13.6
$ws.Columns.Item(3).ColumnWidth = 18.95
$ws.Columns.Item(4).ColumnWidth = 11.6
$ws.Columns.Item(5).ColumnWidth = 15.6
$ws.Columns.Item(6).ColumnWidth = 15.76

# --- Restore selection to F9 ---
[void]$ws.Range("F9").Select()
